$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (rs) for rows 2-18 was 125, should be 63 (root s of 63 GeV)
$ws.Range("F2:F18").Value = 63
